$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.615.16"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.936.35"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'246.26"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'0.9998"

$ws.Range("D7").Value = "'0.4842"
$ws.Range("E7").Value = "  +2.37%  "

$ws.Range("D8").Value = "'0.2920"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "'0.06798"

$ws.Range("D10").Value = "'113.15"
$ws.Range("E10").Value = "  +5.35%  "

$ws.Range("D11").Value = "'19.45"
$ws.Range("E11").Value = "  +4.35%  "

$ws.Range("D12").Value = "1.933.81"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.501"
$ws.Range("E13").Value = "  +2.84%  "

$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.07603"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("D15").Value = "'0.6814"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").Value = "'298.89"
$ws.Range("E16").Value = "  +2.72%  "

$ws.Range("D17").Value = "30.623.32"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "'13.12"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "'0.000007657"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "'0.9998"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "2.189.98"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "'5.562"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").Value = "'0.9993"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'6.532"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").Value = "'9.588"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").Value = "'168.75"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("D27").Value = "'20.37"
$ws.Range("E27").Value = "  -2.53%  "

$ws.Range("D28").Value = "'2.127"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").Value = "'0.1072"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'1.431"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").Value = "'4.198"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").Value = "'4.108"
$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").Value = "'0.05010"
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("D34").Value = "'0.7534"
$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("D35").Value = "'1.151"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").Value = "'0.02044"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("D39").Value = "'2.027"
$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("D40").Value = "'110.23"
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").Value = "'0.4480"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("D42").Value = "'0.8737"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("D43").Value = "'5.831"
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("D44").Value = "'70.57"
$ws.Range("E44").Value = "  +3.54%  "

$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").Value = "'7.347"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "'49.60"
$ws.Range("E47").Value = "  +1.48%  "

$ws.Range("D48").Value = "'9.331"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("D49").Value = "'0.1237"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("D50").Value = "'0.2554"
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").Value = "'35.18"
$ws.Range("E51").Value = "  -0.59%  "
